$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.204.18"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "3.434.71"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Formula = "'407.41"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("D6").Formula = "'128.16"
$ws.Range("E6").Value = "  -3.25%  "

$ws.Range("D7").Formula = "'0.616"
$ws.Range("E7").Value = "  +4.31%  "

$ws.Range("D8").Formula = "'0.999"
$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").Formula = "'0.744"
$ws.Range("E9").Value = "  +11.07%  "

$ws.Range("D10").Formula = "'0.142"
$ws.Range("E10").Value = "  +18.42%  "

$ws.Range("D11").Formula = "'42.47"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Formula = "'0.140"
$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Formula = "'8.65"
$ws.Range("E13").Value = "  +3.55%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Formula = "'20.09"
$ws.Range("E14").Value = "  +2.69%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Formula = "'0.0000191"
$ws.Range("E15").Value = "  +50.80%  "

$ws.Range("D16").Value = "3.438.90"
$ws.Range("E16").Value = "  +2.11%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "62.196.87"
$ws.Range("E17").Value = "  +1.13%  "

$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Formula = "'1.04"
$ws.Range("E18").Value = "  +1.57%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Formula = "'11.38"
$ws.Range("E19").Value = "  +1.63%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Formula = "'363.05"
$ws.Range("E20").Value = "  +19.28%  "

$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Formula = "'87.38"
$ws.Range("E21").Value = "  +5.44%  "

$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Formula = "'3.15"
$ws.Range("E22").Value = "  -1.87%  "

$ws.Range("D23").Formula = "'13.03"
$ws.Range("E23").Value = "  +1.78%  "

$ws.Range("D24").Formula = "'3.17"
$ws.Range("E24").Value = "  +1.28%  "

$ws.Range("D25").Formula = "'31.39"
$ws.Range("E25").Value = "  +7.25%  "

$ws.Range("D26").Formula = "'4.75"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").Formula = "'8.23"
$ws.Range("E27").Value = "  -2.21%  "

$ws.Range("D28").Formula = "'7.70"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Formula = "'2.73"
$ws.Range("E29").Value = "  +9.64%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Formula = "'43.75"
$ws.Range("E30").Value = "  +5.62%  "

$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Formula = "'0.171"
$ws.Range("E31").Value = "  -1.13%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Formula = "'0.115"
$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("D33").Formula = "'11.70"
$ws.Range("E33").Value = "  +3.50%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Formula = "'0.0492"
$ws.Range("E35").Value = "  +2.86%  "

$ws.Range("D36").Formula = "'51.89"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").Formula = "'0.997"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").Formula = "'3.35"
$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("D39").Formula = "'2.91"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Formula = "'143.47"
$ws.Range("E40").Value = "  +4.33%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Formula = "'0.312"
$ws.Range("E41").Value = "  +7.69%  "

$ws.Range("E42").Value = "  +4.48%  "

$ws.Range("D43").Formula = "'1.96"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("D44").Formula = "'3.91"
$ws.Range("E44").Value = "  -1.54%  "

$ws.Range("D45").Formula = "'16.64"
$ws.Range("E45").Value = "  -0.60%  "

$ws.Range("E46").Value = "  +4.44%  "

$ws.Range("D47").Formula = "'21.58"
$ws.Range("E47").Value = "  +1.01%  "

$ws.Range("D48").Value = "2.108.31"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Formula = "'2.30"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Formula = "'1.93"
$ws.Range("E50").Value = "  +2.85%  "

$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Formula = "'0.0366"
$ws.Range("E51").Value = "  +9.08%  "
